$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1756.6666
$ws.Range("J40").Value = 2069.8572
$ws.Range("L40").Value = 2069.8572
$ws.Range("N40").Value = -2419.8572
$ws.Range("H43").Value = 2045.0834
$ws.Range("I43").Value = 1250
$ws.Range("J43").Value = 2613
$ws.Range("K43").Value = 1250
$ws.Range("L43").Value = 2613
$ws.Range("M43").Value = -1181
$ws.Range("N43").Value = -2751
$ws.Range("H64").Value = 3614.3333
$ws.Range("I64").Value = 3818.6667
$ws.Range("K64").Value = 3818.6667
$ws.Range("M64").Value = -3570.6667
$ws.Range("H67").Value = 3614.3333
$ws.Range("I67").Value = 3818.6667
$ws.Range("K67").Value = 3818.6667
$ws.Range("M67").Value = -2960.6667
$ws.Range("H76").Value = 2718.889
$ws.Range("I76").Value = 2660.8333
$ws.Range("J76").Value = 2835
$ws.Range("K76").Value = 2660.8333
$ws.Range("L76").Value = 2835
$ws.Range("M76").Value = -2345.8333
$ws.Range("N76").Value = -3465
$ws.Range("H79").Value = 2718.889
$ws.Range("I79").Value = 2660.8333
$ws.Range("J79").Value = 2835
$ws.Range("K79").Value = 2660.8333
$ws.Range("L79").Value = 2835
$ws.Range("M79").Value = -1568.8333
$ws.Range("N79").Value = -5019
$ws.Range("H80").Value = 11667.167
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 16500.75
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 49502.25
$ws.Range("M80").Value = -5002
$ws.Range("N80").Value = -51498.25
$ws.Range("H83").Value = 11667.167
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 16500.75
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 148506.75
$ws.Range("M83").Value = -13008
$ws.Range("N83").Value = -158490.75
$ws.Range("H134").Value = 25210.834
$ws.Range("J134").Value = 25210.834
$ws.Range("L134").Value = 25210.834
$ws.Range("N134").Value = -35350.834

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15576.161
$ws.Range("I32").Value = 13984.551
$ws.Range("J32").Value = 21575.309
$ws.Range("K32").Value = 13984.551
$ws.Range("L32").Value = 21575.309
$ws.Range("M32").Value = -13697.551
$ws.Range("N32").Value = -22149.309
$ws.Range("H33").Value = 18432.143
$ws.Range("I33").Value = 17012.5
$ws.Range("K33").Value = 17012.5
$ws.Range("M33").Value = -16683.5
$ws.Range("H113").Value = 30000
$ws.Range("J113").Value = 30000
$ws.Range("L113").Value = 30000
$ws.Range("N113").Value = -38678
$ws.Range("H122").Value = 2134.342
$ws.Range("I122").Value = 1589.2667
$ws.Range("J122").Value = 4178.375
$ws.Range("K122").Value = 4767.800099999999
$ws.Range("L122").Value = 12535.125
$ws.Range("M122").Value = -2317.800099999999
$ws.Range("N122").Value = -17435.125
$ws.Range("H134").Value = 31157.143
$ws.Range("J134").Value = 31157.143
$ws.Range("L134").Value = 31157.143
$ws.Range("N134").Value = -41297.143

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 29828.889
$ws.Range("J51").Value = 29828.889
$ws.Range("L51").Value = 29828.889
$ws.Range("N51").Value = -30810.889
$ws.Range("H105").Value = 1646.7894
$ws.Range("I105").Value = 1329.875
$ws.Range("J105").Value = 3337
$ws.Range("K105").Value = 1329.875
$ws.Range("L105").Value = 3337
$ws.Range("M105").Value = 417.125
$ws.Range("N105").Value = -6831

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 100002
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 100002
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 100002
$ws.Range("N4").Value = -100226
$ws.Range("H22").Value = 971.1667
$ws.Range("I22").Value = 337
$ws.Range("K22").Value = 337
$ws.Range("M22").Value = 13
$ws.Range("H31").Value = 4069.6072
$ws.Range("I31").Value = 3104.96
$ws.Range("J31").Value = 4847.5483
$ws.Range("K31").Value = 3104.96
$ws.Range("L31").Value = 4847.5483
$ws.Range("M31").Value = -2809.96
$ws.Range("N31").Value = -5437.5483
$ws.Range("H34").Value = 4069.6072
$ws.Range("I34").Value = 3104.96
$ws.Range("J34").Value = 4847.5483
$ws.Range("K34").Value = 3104.96
$ws.Range("L34").Value = 4847.5483
$ws.Range("M34").Value = -2902.96
$ws.Range("N34").Value = -5251.5483
$ws.Range("H62").Value = 3780
$ws.Range("I62").Value = 2431.5625
$ws.Range("J62").Value = 6476.875
$ws.Range("K62").Value = 2431.5625
$ws.Range("L62").Value = 6476.875
$ws.Range("M62").Value = -1807.5625
$ws.Range("N62").Value = -7724.875
$ws.Range("H65").Value = 3780
$ws.Range("I65").Value = 2431.5625
$ws.Range("J65").Value = 6476.875
$ws.Range("K65").Value = 12157.8125
$ws.Range("L65").Value = 32384.375
$ws.Range("M65").Value = -9037.8125
$ws.Range("N65").Value = -38624.375
$ws.Range("M4").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 542.75
$ws.Range("J4").Value = 700
$ws.Range("L4").Value = 2100
$ws.Range("N4").Value = -2324
$ws.Range("H12").Value = 422.5
$ws.Range("J12").Value = 639.2273
$ws.Range("L12").Value = 1917.6819
$ws.Range("N12").Value = -2263.6819
$ws.Range("H92").Value = 1537
$ws.Range("J92").Value = 1486.8667
$ws.Range("L92").Value = 4460.6001
$ws.Range("N92").Value = -6956.6001
$ws.Range("H97").Value = 763.44446
$ws.Range("I97").Value = 422.5
$ws.Range("J97").Value = 1036.2
$ws.Range("K97").Value = 1267.5
$ws.Range("L97").Value = 3108.6
$ws.Range("M97").Value = -771.5
$ws.Range("N97").Value = -4100.6
$ws.Range("H98").Value = 241.27777
$ws.Range("I98").Value = 200.23077
$ws.Range("J98").Value = 348
$ws.Range("K98").Value = 600.69231
$ws.Range("L98").Value = 1044
$ws.Range("M98").Value = 897.30769
$ws.Range("N98").Value = -4040
$ws.Range("H107").Value = 1534
$ws.Range("J107").Value = 2182.6667
$ws.Range("L107").Value = 6548.000100000001
$ws.Range("N107").Value = -10388.0001
$ws.Range("H131").Value = 1674.9736
$ws.Range("I131").Value = 3717
$ws.Range("J131").Value = 1213.871
$ws.Range("K131").Value = 11151
$ws.Range("L131").Value = 3641.613
$ws.Range("M131").Value = -6111
$ws.Range("N131").Value = -13721.613
$ws.Range("H133").Value = 6097.5
$ws.Range("I133").Value = 6335.364
$ws.Range("J133").Value = 5723.7144
$ws.Range("K133").Value = 19006.092
$ws.Range("L133").Value = 17171.1432
$ws.Range("M133").Value = -13946.092
$ws.Range("N133").Value = -27291.1432

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 10823
$ws.Range("J46").Value = 11430.667
$ws.Range("L46").Value = 11430.667
$ws.Range("N46").Value = -11742.667
$ws.Range("H70").Value = 4915
$ws.Range("I70").Value = 5071.4287
$ws.Range("J70").Value = 4550
$ws.Range("K70").Value = 5071.4287
$ws.Range("L70").Value = 4550
$ws.Range("M70").Value = -4801.4287
$ws.Range("N70").Value = -5090
$ws.Range("H73").Value = 4915
$ws.Range("I73").Value = 5071.4287
$ws.Range("J73").Value = 4550
$ws.Range("K73").Value = 5071.4287
$ws.Range("L73").Value = 4550
$ws.Range("M73").Value = -4135.4287
$ws.Range("N73").Value = -6422
$ws.Range("H80").Value = 2950
$ws.Range("I80").Value = 3153.8462
$ws.Range("J80").Value = 2709.0908
$ws.Range("K80").Value = 3153.8462
$ws.Range("L80").Value = 2709.0908
$ws.Range("M80").Value = -2155.8462
$ws.Range("N80").Value = -4705.0908
$ws.Range("H83").Value = 2950
$ws.Range("I83").Value = 3153.8462
$ws.Range("J83").Value = 2709.0908
$ws.Range("K83").Value = 15769.231
$ws.Range("L83").Value = 13545.454
$ws.Range("M83").Value = -10777.231
$ws.Range("N83").Value = -23529.454

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 21007.363
$ws.Range("J17").Value = 23008.1
$ws.Range("L17").Value = 23008.1
$ws.Range("N17").Value = -23348.1
$ws.Range("H35").Value = 8896.5
$ws.Range("I35").Value = 234.4
$ws.Range("J35").Value = 23333.334
$ws.Range("K35").Value = 234.4
$ws.Range("L35").Value = 23333.334
$ws.Range("M35").Value = 101.6
$ws.Range("N35").Value = -24005.334
$ws.Range("H46").Value = 3142.5715
$ws.Range("I46").Value = 399.6
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 399.6
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = -211.6
$ws.Range("N46").Value = -10376
$ws.Range("H106").Value = 29500
$ws.Range("J106").Value = 29500
$ws.Range("L106").Value = 29500
$ws.Range("N106").Value = -32024
$ws.Range("H122").Value = 3266.75
$ws.Range("I122").Value = 2648.32
$ws.Range("J122").Value = 4672.273
$ws.Range("K122").Value = 7944.960000000001
$ws.Range("L122").Value = 14016.819
$ws.Range("M122").Value = -5494.960000000001
$ws.Range("N122").Value = -18916.819

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 100006
$ws.Range("J13").Value = 100006
$ws.Range("L13").Value = 100006
$ws.Range("N13").Value = -100286
$ws.Range("H55").Value = 7197.8
$ws.Range("J55").Value = 8247.25
$ws.Range("L55").Value = 8247.25
$ws.Range("N55").Value = -8801.25
$ws.Range("H81").Value = 2375
$ws.Range("H84").Value = 2375
$ws.Range("H126").Value = 36566.7
$ws.Range("I126").Value = 49862.57
$ws.Range("J126").Value = 5543
$ws.Range("K126").Value = 149587.71
$ws.Range("L126").Value = 16629
$ws.Range("M126").Value = -147117.71
$ws.Range("N126").Value = -21569
